$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A8").Value = 10
$ws.Range("B8").Value = "number_of_seasons"
$ws.Range("C8").Value = 0.6920467605659095
$ws.Range("D8").Value = 0.1822935422608884
